$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 data (appended test case row).
# A12 looks like a pure number ("141287"); format the cell as Text first so
# it is stored as a text value (matching t="str" in the source XML) instead
# of being coerced into a numeric cell.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "141287"
$ws.Range("B12").Value = "Client->Client Summary->Training->To verify that when the user clicks the Save button on the Client Record-Training Data screen, the data should be saved, and they should be navigated to the training main screen."
$ws.Range("C12").Value = "assert,click,assert,click"
$ws.Range("D12").Value = "Trainings,no value,fixture,no value"
$ws.Range("E12").Value = "training_assert,click_training,first_name,finish"
$ws.Range("F12").Value = "yes"
